$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Student 1 (row 6): KISW score updated, totals/average recalculated
$ws.Range("E6").Value = 78
$ws.Range("K6").Value = 545
$ws.Range("L6").Value = 68.125

# Student 2 (row 7): KISW score updated, totals/average recalculated
$ws.Range("E7").Value = 60
$ws.Range("K7").Value = 541
$ws.Range("L7").Value = 67.625

# SubjectTotal row (row 8): KISW column total updated
$ws.Range("E8").Value = 138
$ws.Range("K8").Value = 1086
$ws.Range("L8").Value = 135.75

# SubjectAverage row (row 9): KISW column average updated
$ws.Range("E9").Value = 69
$ws.Range("K9").Value = 543
$ws.Range("L9").Value = 67.875

# SubjectGrades row (row 10): KISW grade improved from C- to B-
$ws.Range("E10").Value = "B-"
